$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand symbol / Receptor symbol / Target cluster, then the
# 16 numeric metric columns (E..T), for rows 2..7 of the LR-pair table.
$rows = @(
    ,@("FAPs","Angpt4","Tek","ECs", 3,1,0.6486246666666667,1.945874,0.7736027733825248,0.7736027733825248,3,1,39.413899,118.241697,0.7310109930444597,0.7310109930444597,25.56482710090867,230.083443908178,0.5655121315923076,0.5655121315923076)
    ,@("FAPs","Angpt4","Tek","FAPs", 3,1,0.6486246666666667,1.945874,0.7736027733825248,0.7736027733825248,3,1,11.40791933333333,34.223758,0.2115830874897988,0.2115830874897987,7.399457874943557,66.59512087449201,0.1636812632829457,0.1636812632829457)
    ,@("FAPs","Angpt4","Tek","sCs", 3,1,0.6486246666666667,1.945874,0.7736027733825248,0.7736027733825248,3,1,3.095153333333334,9.28546,0.05740591946574151,0.0574059194657415,2.007592799115556,18.06833519204,0.0444093785072715,0.04440937850727149)
    ,@("sCs","Angpt4","Tek","ECs", 2,0.6666666666666666,0.189822,0.569466,0.2263972266174752,0.2263972266174752,3,1,39.413899,118.241697,0.7310109930444597,0.7310109930444597,7.481625135978001,67.334626223802,0.1654988614521521,0.1654988614521521)
    ,@("sCs","Angpt4","Tek","FAPs", 2,0.6666666666666666,0.189822,0.569466,0.2263972266174752,0.2263972266174752,3,1,11.40791933333333,34.223758,0.2115830874897988,0.2115830874897987,2.165474063692,19.489266573228,0.04790182420685304,0.04790182420685304)
    ,@("sCs","Angpt4","Tek","sCs", 2,0.6666666666666666,0.189822,0.569466,0.2263972266174752,0.2263972266174752,3,1,3.095153333333334,9.28546,0.05740591946574151,0.0574059194657415,0.5875281960400001,5.287753764360001,0.01299654095847001,0.01299654095847001)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
